# Apply the commit's changes to the StructureDefinition-valid-age-reason workbook.
# Sheet 1 = "Metadata" (Property/Value table)
# Sheet 2 = "Elements" (big element definition table)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 (Metadata) changes ---

# 1. URL value: pythia -> cicada
$ws1.Cells.Item(2, 2).Value = "http://fhirfli.dev/fhir/ig/cicada/StructureDefinition/valid-age-reason"

# 2. Date value updated
$ws1.Cells.Item(8, 2).Value = "2026-02-11T14:37:07-05:00"

# 3. Insert a new "Jurisdiction" row right after "Contact" (row 10), pushing
#    Description/Purpose/Copyright/etc down by one row.
$ws1.Rows.Item(11).Insert()

# Copy formatting (style) from the row above so the new row matches the rest
# of the table (border/alignment), then fill in the new values.
$ws1.Range("A10:B10").Copy()
$ws1.Range("A11:B11").PasteSpecial(-4122)
$ws1.Cells.Item(11, 1).Value = "Jurisdiction"
$ws1.Cells.Item(11, 2).Value = ""

# --- Sheet2 (Elements) changes ---

# Same pythia -> cicada URL swaps (StructureDefinition url + ValueSet url)
$ws2.Cells.Item(5, 18).Value = "http://fhirfli.dev/fhir/ig/cicada/StructureDefinition/valid-age-reason"
$ws2.Cells.Item(6, 26).Value = "http://fhirfli.dev/fhir/ig/cicada/ValueSet/valid-age-reason"
